# Generate Report for Handback
# Update the "generate date" / handoff-handback timestamp cells in the
# Overview, zh-cn and de-de sheets to reflect a freshly regenerated report.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview sheet: Latest HO Xliff Generate Date (also mirrored on de-de!H2)
$wsOverview.Range("G2").Value = "2016-08-27 23:03:19"

# zh-cn sheet: Correspond Handoff Datetime / Correspond Handback DateTime
$wsZhCn.Range("H2").Value = "2016-08-27 23:03:14"
$wsZhCn.Range("K2").Value = "2016-08-27 23:03:44"

# de-de sheet: Correspond Handoff Datetime (mirrors Overview!G2) / Correspond Handback DateTime
$wsDeDe.Range("H2").Value = "2016-08-27 23:03:19"
$wsDeDe.Range("K2").Value = "2016-08-27 23:03:51"
